$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Row 6 holds the "Other" category keyword list, stored alphabetically
# across the row starting at column B (A6 is the category label
# "Other"). Read the existing keywords, insert the new keyword
# "Binding waste" in its correct alphabetical spot (between "Biblical"
# and "Cartography"), and write the whole row back in one shot so nothing
# in other rows is disturbed.

$lastCol = $ws.Cells.Item(6, $ws.Columns.Count).End(-4159).Column  # xlToLeft
$existing = @()
for ($c = 2; $c -le $lastCol; $c++) {
    $existing += ,$ws.Cells.Item(6, $c).Value()
}

$newWord = "Binding waste"
$insertAt = $existing.Count
for ($i = 0; $i -lt $existing.Count; $i++) {
    if ($newWord -lt $existing[$i]) {
        $insertAt = $i
        break
    }
}

$updated = @()
if ($insertAt -ge 1) {
    $updated += $existing[0..($insertAt - 1)]
}
$updated += $newWord
if ($insertAt -le $existing.Count - 1) {
    $updated += $existing[$insertAt..($existing.Count - 1)]
}

for ($i = 0; $i -lt $updated.Count; $i++) {
    $ws.Cells.Item(6, $i + 2).Value = $updated[$i]
}

# Make sure every keyword cell (including the newly created one at the
# end of the row) keeps the same formatting used throughout the rest of
# the row.
$styleSource = $ws.Cells.Item(6, 2)
$destRange = $ws.Range($ws.Cells.Item(6, 2), $ws.Cells.Item(6, $updated.Count + 1))
$styleSource.Copy() | Out-Null
$destRange.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
